$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.586.65'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.509.70'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.72%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.84%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.78'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.82%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.586'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.25%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.539'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.02'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0814'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.87%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.61'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.53%  '
$ws.Range("E13").Value = '  -2.45%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.900.96'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.54'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +7.17%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.506.75'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.69%  '
$ws.Range("E17").Value = '  -2.48%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.540.66'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.24%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.91'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.09%  '
$ws.Range("E20").Value = '  -1.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.49'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.56'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '253.23'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.77%  '
$ws.Range("E24").Value = '  +1.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.03'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.91%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '27.02'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.36%  '
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("E28").Value = '  +12.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.13'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '37.85'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.92'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.68%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '154.84'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.60%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.22'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.33'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.07'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.91%  '
$ws.Range("E36").Value = '  -1.70%  '
$ws.Range("E37").Value = '  -4.55%  '
$ws.Range("E38").Value = '  -0.82%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '24.37'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.04%  '
$ws.Range("E40").Value = '  +1.17%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.86'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.42%  '
$ws.Range("E42").Value = '  +0.46%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.02'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.18%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("E45").Value = '  -0.06%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.023.85'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '84.46'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.68%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.94'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.756.64'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.62%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '73.26'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.190'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.56%  '
